$ws = $excel.ActiveWorkbook.ActiveSheet

# Force text/string storage for D-column cells whose new numeric-looking
# values would otherwise be auto-coerced to the Number type by Excel,
# since the source data keeps these as plain text cells.
$textCells = @("D41", "D29", "D32", "D45", "D36", "D34", "D23", "D27", "D22", "D24", "D49", "D6", "D11", "D5", "D13", "D19", "D25", "D30", "D37", "D39", "D47", "D31", "D28", "D20", "D42", "D10", "D14", "D40", "D21", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated price / volume figures scraped for this run.
$ws.Range("D2").Value = "64.345.23"
$ws.Range("E2").Value = "  -3.50%  "
$ws.Range("D3").Value = "3.162.00"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "608.43"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "147.97"
$ws.Range("E6").Value = "  -7.09%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.161.91"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").Value = "5.55"
$ws.Range("E11").Value = "  -6.56%  "
$ws.Range("E12").Value = "  -6.13%  "
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  -7.25%  "
$ws.Range("D14").Value = "35.84"
$ws.Range("E14").Value = "  -9.37%  "
$ws.Range("D15").Value = "3.679.48"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "64.340.01"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "3.158.03"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "6.97"
$ws.Range("E19").Value = "  -5.79%  "
$ws.Range("D20").Value = "483.32"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("D21").Value = "14.88"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  -5.83%  "
$ws.Range("D23").Value = "7.77"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").Value = "13.77"
$ws.Range("E24").Value = "  -7.39%  "
$ws.Range("D25").Value = "83.88"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").Value = "8.52"
$ws.Range("E28").Value = "  -6.45%  "
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -8.11%  "
$ws.Range("D30").Value = "6.82"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "0.114"
$ws.Range("E31").Value = "  -29.49%  "
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "26.31"
$ws.Range("E34").Value = "  -6.70%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  -6.67%  "
$ws.Range("D37").Value = "54.35"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("E38").Value = "  -9.74%  "
$ws.Range("D39").Value = "456.13"
$ws.Range("E39").Value = "  -8.17%  "
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -13.58%  "
$ws.Range("D41").Value = "0.0397"
$ws.Range("E41").Value = "  -7.49%  "
$ws.Range("D42").Value = "8.47"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("E43").Value = "  -8.14%  "
$ws.Range("D44").Value = "2.853.65"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  -9.47%  "
$ws.Range("E46").Value = "  -8.48%  "
$ws.Range("D47").Value = "26.51"
$ws.Range("E47").Value = "  -7.64%  "
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("E50").Value = "  -4.50%  "
$ws.Range("D51").Value = "118.88"
$ws.Range("E51").Value = "  -2.34%  "
